$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SpecsDataCalib")
for ($c = 1; $c -le 10; $c++) {
    $ws.Columns.Item($c).AutoFit()
}
